$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-24 Monday" "2024-06-25 Tuesday"

Replace-Text "906×4=3624" "871×5=4355"
Replace-Text "345×9=3105" "682×9=6138"
Replace-Text "521×9=4689" "732×5=3660"
Replace-Text "867×2=1734" "951×2=1902"
Replace-Text "187×9=1683" "806×5=4030"

Replace-Text "907×8=7256" "795×5=3975"
Replace-Text "597×9=5373" "755×5=3775"
Replace-Text "293×7=2051" "222×6=1332"
Replace-Text "352×5=1760" "577×4=2308"
Replace-Text "447×3=1341" "910×5=4550"

Replace-Text "620×2=1240" "805×2=1610"
Replace-Text "388×4=1552" "434×8=3472"
Replace-Text "847×4=3388" "744×7=5208"
Replace-Text "867×5=4335" "130×6=780"
Replace-Text "214×8=1712" "852×8=6816"

Replace-Text "354×9=3186" "880×3=2640"
Replace-Text "116×5=580" "710×4=2840"
Replace-Text "352×6=2112" "947×3=2841"
Replace-Text "694×2=1388" "387×3=1161"
Replace-Text "983×2=1966" "639×6=3834"

Replace-Text "989×6=5934" "754×7=5278"
Replace-Text "450×5=2250" "284×8=2272"
Replace-Text "175×5=875" "826×3=2478"
Replace-Text "431×3=1293" "932×7=6524"
Replace-Text "419×6=2514" "596×6=3576"
